$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, reusing the same style as the other
# header cells (e.g. G1) so the new column matches the existing formatting.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new column's data values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
